# Revert the wide 6-column "one row of labels / one row of values" layout
# back to a tall 2-column "parameter / value" table (one parameter name +
# its value per row), with bold headers and a text-formatted value column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writes a plain string into a cell without letting Excel's input parser
# reinterpret text that looks like a boolean/number/date (e.g. "TRUE",
# "FALSE") as anything other than literal text: build it as a quoted
# string formula, then convert that formula to a static value in place
# (so the cell ends up as a plain shared-string cell, not a formula).
function Set-TextValue($cell, [string]$text) {
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = $false
}

# Capture the existing wide-format data (row 1 = parameter names,
# row 2 = values, six parameter/value pairs spread across columns A:F)
# before touching anything.
$pairs = @()
for ($col = 1; $col -le 6; $col++) {
    $paramName = $ws.Cells.Item(1, $col).Value()
    $paramValue = $ws.Cells.Item(2, $col).Value()
    $pairs += , @($paramName, $paramValue)
}

# Wipe the sheet and any leftover column formatting from the old
# 6-column layout so nothing lingers into the new 2-column one.
$ws.Cells.Clear()
$ws.Columns.Item(1).ClearFormats()
$ws.Columns.Item(2).ClearFormats()

# Header row: bold "parameter" / "value" labels.
Set-TextValue $ws.Cells.Item(1, 1) "parameter"
Set-TextValue $ws.Cells.Item(1, 2) "value"
$ws.Range("A1:B1").Font.Bold = $true
$ws.Cells.Item(1, 2).NumberFormat = "@"

# Data rows: one parameter/value pair per row, starting at row 2.
$row = 2
foreach ($pair in $pairs) {
    Set-TextValue $ws.Cells.Item($row, 1) $pair[0]

    $valueCell = $ws.Cells.Item($row, 2)
    $valueCell.NumberFormat = "@"
    if ($pair[1] -is [bool]) {
        $valueCell.Value = $pair[1]
    } else {
        Set-TextValue $valueCell $pair[1]
    }
    $row++
}

# Column A sizes itself to its (now longer) parameter-name content;
# column B keeps the sheet's original default width.
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).ColumnWidth = 9.140625

$ws.Range("B6").Select()
